$d = $word.ActiveDocument

# 1.7.1.1
$p = $d.Paragraphs.Item(2)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3003BB05" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.1.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 8 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 42 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 58</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.1.2
$p = $d.Paragraphs.Item(3)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1A2A716E" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.1.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 47 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 53</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.1.3
$p = $d.Paragraphs.Item(4)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="694E4ECB" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.1.3 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 44 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 56</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.1.4
$p = $d.Paragraphs.Item(5)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7C31A6AC" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.1.4 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 9 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 41 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 59</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.1.5
$p = $d.Paragraphs.Item(6)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="17D40807" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.1.5 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 47 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 53</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.1.6
$p = $d.Paragraphs.Item(7)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3C44AB9E" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.1.6 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 65 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 69 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 72</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.2.1
$p = $d.Paragraphs.Item(8)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="31789348" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.2.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 45 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 55</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.2.2
$p = $d.Paragraphs.Item(9)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="79B7048C" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.2.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 8 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 42 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 58</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.2.3
$p = $d.Paragraphs.Item(10)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7D36756D" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.2.3 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 9 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 41 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 59</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.2.4
$p = $d.Paragraphs.Item(11)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2D43DF7E" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.2.4 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 10 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 7 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 59 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 70 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 80</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.3.1
$p = $d.Paragraphs.Item(12)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1FE1B938" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.3.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 8 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 40 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 60</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.3.2
$p = $d.Paragraphs.Item(13)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="318B6805" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.3.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 21 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 29 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 71</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.3.3
$p = $d.Paragraphs.Item(14)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4C225704" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.3.3 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 9 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 41 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 59</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.3.4
$p = $d.Paragraphs.Item(15)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5C710E11" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.3.4 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 9 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 55 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 65 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 74</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.4.1
$p = $d.Paragraphs.Item(16)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7835854E" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.4.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 12 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 38 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 62</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.4.2
$p = $d.Paragraphs.Item(17)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3847E717" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.4.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 8 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 7 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 42 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 58</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.4.3
$p = $d.Paragraphs.Item(18)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="39792A9B" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.4.3 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 45 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 55</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.4.4
$p = $d.Paragraphs.Item(19)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7864608F" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.4.4 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 11 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 39 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 61</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.4.5
$p = $d.Paragraphs.Item(20)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3DD94F80" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.4.5 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 44 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 56</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.4.6
$p = $d.Paragraphs.Item(21)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="54CBBD96" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.4.6 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 52 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 58 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 63</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.5.1
$p = $d.Paragraphs.Item(22)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="27350666" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.5.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 7 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 43 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 57</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.5.2
$p = $d.Paragraphs.Item(23)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="384347E7" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.5.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 9 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 41 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 59</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.5.3
$p = $d.Paragraphs.Item(24)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="082A4263" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.5.3 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 10 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 40 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 60</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.5.4
$p = $d.Paragraphs.Item(25)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="58417ECB" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.5.4 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 13 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 43 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 57 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 70</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.6.1
$p = $d.Paragraphs.Item(26)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="63AA1FBD" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.6.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 44 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 56</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.6.2
$p = $d.Paragraphs.Item(27)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="10908861" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.6.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 11 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 39 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 61</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.6.3
$p = $d.Paragraphs.Item(28)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="13CCA83D" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.6.3 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 12 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 37 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 63</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.6.4
$p = $d.Paragraphs.Item(29)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2D5EEC13" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.6.4 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 10 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 40 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 60</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.6.5
$p = $d.Paragraphs.Item(30)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2D936C6C" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.6.5 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 45 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 55</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.6.6
$p = $d.Paragraphs.Item(31)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0BB324B5" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.6.6 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 11 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 39 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 61</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.6.7
$p = $d.Paragraphs.Item(32)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7293F371" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.6.7 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 10 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 8 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 63 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 74 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 84</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.7.1
$p = $d.Paragraphs.Item(33)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="30EF1F07" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.7.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 39 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 56</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.7.2
$p = $d.Paragraphs.Item(34)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="60861573" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.7.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 13 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 48 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 65 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 77</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.8.1
$p = $d.Paragraphs.Item(35)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="12F1C8AD" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.8.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 8 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 40 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 58</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.8.2
$p = $d.Paragraphs.Item(36)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="60AC5DE7" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.8.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 12 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 35 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 59</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.8.3
$p = $d.Paragraphs.Item(37)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4F9B4494" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.8.3 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 8 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 7 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 38 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 56</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.8.4
$p = $d.Paragraphs.Item(38)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="12A2BD8E" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.8.4 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 35 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 42 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 47</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.9.1
$p = $d.Paragraphs.Item(39)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="084EB79F" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.9.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 43 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 56</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.9.2
$p = $d.Paragraphs.Item(40)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4EE709BB" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.9.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 8 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 59 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 67</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.10.1
$p = $d.Paragraphs.Item(41)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="106D22E7" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.10.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 12 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 36 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 59</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.10.2
$p = $d.Paragraphs.Item(42)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="45E0B4D6" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.10.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 9 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 41 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 53</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.10.3
$p = $d.Paragraphs.Item(43)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="58864311" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:lastRenderedPageBreak/><w:t>1.7.10.3 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 28 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 36 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 40</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.11.1
$p = $d.Paragraphs.Item(44)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="31C7CE66" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.11.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 13 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 9 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 37 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 63</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.11.2
$p = $d.Paragraphs.Item(45)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5945F61F" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.11.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 14 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 8 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 31 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 46 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 60</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.12.1
$p = $d.Paragraphs.Item(46)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6BFFDFA5" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.12.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 14 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 35 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 63</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.12.2
$p = $d.Paragraphs.Item(47)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5233BBB8" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.12.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 10 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 7 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 48 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 63 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 71</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.13.1
$p = $d.Paragraphs.Item(48)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="25B56893" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.13.1 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 11 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 42 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 53</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.13.2
$p = $d.Paragraphs.Item(49)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3B798401" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.13.2 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 40 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 54</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.13.3
$p = $d.Paragraphs.Item(50)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3FC9E178" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.13.3 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 4 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 42 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 54</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')

# 1.7.13.4
$p = $d.Paragraphs.Item(51)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3EDBA0E5" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.13.4 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 11 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 6 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 5 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 1 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 35 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 50 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 59</w:t></w:r><w:r><w:t xml:space="preserve"> Ok</w:t></w:r></w:p>')

# 1.7.13.5
$p = $d.Paragraphs.Item(52)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0AB8450A" w14:textId="77777777" w:rsidR="009A36CB" w:rsidRDefault="009A36CB" w:rsidP="009A36CB"><w:r><w:t>1.7.13.5 :</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 10 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 2 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 9 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 3 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 0 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 33 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 51 </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> 60</w:t></w:r><w:r><w:t xml:space="preserve"> OK</w:t></w:r></w:p>')
